$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Additional/updated data rows (columns A:T) for rows 2 through 7.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric expression / specificity metrics.
$data = @(
    @("FAPs", "Ccl21b", "Ackr4", "ECs",  3, 1, 0.3964663333333333, 1.189399,            0.8220900069740165, 0.8220900069740164, 1, 0.3333333333333333, 0.1143813333333333, 0.343144,            0.03414257747883775, 0.03414257747883775, 0.04534834782844444, 0.408135130456,      0.02806827175768862, 0.02806827175768862),
    @("FAPs", "Ccl21b", "Ackr4", "FAPs", 3, 1, 0.3964663333333333, 1.189399,            0.8220900069740165, 0.8220900069740164, 3, 1,                   3.175982333333334,  9.527947000000001,  0.9480237703755849,  0.9480237703755849,  1.259170070428111,   11.332530633853,     0.779360867999598,   0.7793608679995979),
    @("FAPs", "Ccl21b", "Ackr4", "sCs",  3, 1, 0.3964663333333333, 1.189399,            0.8220900069740165, 0.8220900069740164, 2, 0.6666666666666666, 0.05974466666666667, 0.179234,            0.01783365214557738, 0.01783365214557738, 0.02368674892955555, 0.213180740366,      0.01466086721672989, 0.01466086721672989),
    @("sCs",  "Ccl21b", "Ackr4", "ECs",  1, 0.3333333333333333, 0.0858,             0.2574,              0.1779099930259836, 0.1779099930259836, 1, 0.3333333333333333, 0.1143813333333333, 0.343144,            0.03414257747883775, 0.03414257747883775, 0.0098139184,         0.08832526560000001, 0.006074305721149128, 0.006074305721149128),
    @("sCs",  "Ccl21b", "Ackr4", "FAPs", 1, 0.3333333333333333, 0.0858,             0.2574,              0.1779099930259836, 0.1779099930259836, 3, 1,                   3.175982333333334,  9.527947000000001,  0.9480237703755849,  0.9480237703755849,  0.2724992842,         2.4524935578,        0.168662902375987,   0.168662902375987),
    @("sCs",  "Ccl21b", "Ackr4", "sCs",  1, 0.3333333333333333, 0.0858,             0.2574,              0.1779099930259836, 0.1779099930259836, 2, 0.6666666666666666, 0.05974466666666667, 0.179234,            0.01783365214557738, 0.01783365214557738, 0.005126092400000001, 0.0461348316,        0.003172784928847489, 0.003172784928847489)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowVals[$c]
    }
}
